$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" values (column E) for rows 17-22
$ws.Range("E17").Value = "2206"
$ws.Range("E18").Value = "2207"
$ws.Range("E19").Value = "2208"
$ws.Range("E20").Value = "2508"
$ws.Range("E21").Value = "2508"
$ws.Range("E22").Value = "2508"

# Swap "Valor Mora" values (column F) between rows 17 and 20
$ws.Range("F17").Value = 40000
$ws.Range("F20").Value = 56940
